# Automatic update of files.
# Updates the "Förändrad" (changed) date stamp to 46077 for all rows,
# and refreshes the Hammarö logging-notification rows (9-33) with the
# latest data pulled from the source (new/removed entries cause the
# existing rows to shift and be overwritten in place).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 46077
# Row 3
$ws.Range("C3").Value = 46077
# Row 4
$ws.Range("C4").Value = 46077
# Row 5
$ws.Range("C5").Value = 46077
# Row 6
$ws.Range("C6").Value = 46077
# Row 7
$ws.Range("C7").Value = 46077
# Row 8
$ws.Range("C8").Value = 46077
# Row 9
$ws.Range("A9").Value = 'A 58382-2024'
$ws.Range("B9").Value = 45632
$ws.Range("C9").Value = 46077
$ws.Range("G9").Value = 1.6
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1
$ws.Range("P9").Value = 1
$ws.Range("R9").Value = 'Knärot'
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 58382-2024 artfynd.xlsx", "A 58382-2024")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 58382-2024 karta.png", "A 58382-2024")'
$ws.Range("U9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/knärot/A 58382-2024 karta knärot.png", "A 58382-2024")'
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 58382-2024 FSC-klagomål.docx", "A 58382-2024")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 58382-2024 FSC-klagomål mail.docx", "A 58382-2024")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 58382-2024 tillsynsbegäran.docx", "A 58382-2024")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 58382-2024 tillsynsbegäran mail.docx", "A 58382-2024")'
$ws.Range("Z9").ClearContents()
# Row 10
$ws.Range("A10").Value = 'A 55068-2023'
$ws.Range("B10").Value = 45237
$ws.Range("C10").Value = 46077
$ws.Range("F10").Value = 'Kommuner'
$ws.Range("G10").Value = 2.4
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("R10").Value = 'Motaggsvamp'
$ws.Range("S10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 55068-2023 artfynd.xlsx", "A 55068-2023")'
$ws.Range("T10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 55068-2023 karta.png", "A 55068-2023")'
$ws.Range("U10").ClearContents()
$ws.Range("V10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 55068-2023 FSC-klagomål.docx", "A 55068-2023")'
$ws.Range("W10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 55068-2023 FSC-klagomål mail.docx", "A 55068-2023")'
$ws.Range("X10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 55068-2023 tillsynsbegäran.docx", "A 55068-2023")'
$ws.Range("Y10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 55068-2023 tillsynsbegäran mail.docx", "A 55068-2023")'
# Row 11
$ws.Range("A11").Value = 'A 20755-2021'
$ws.Range("B11").Value = 44316
$ws.Range("C11").Value = 46077
$ws.Range("F11").ClearContents()
$ws.Range("G11").Value = 16.8
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("R11").Value = 'Gulsparv'
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 20755-2021 artfynd.xlsx", "A 20755-2021")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 20755-2021 karta.png", "A 20755-2021")'
$ws.Range("U11").ClearContents()
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 20755-2021 FSC-klagomål.docx", "A 20755-2021")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 20755-2021 FSC-klagomål mail.docx", "A 20755-2021")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 20755-2021 tillsynsbegäran.docx", "A 20755-2021")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 20755-2021 tillsynsbegäran mail.docx", "A 20755-2021")'
$ws.Range("Z11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/fåglar/A 20755-2021 prioriterade fågelarter.docx", "A 20755-2021")'
# Row 12
$ws.Range("A12").Value = 'A 24618-2022'
$ws.Range("B12").Value = 44727
$ws.Range("C12").Value = 46077
$ws.Range("G12").Value = 4.6
$ws.Range("H12").Value = 1
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1
$ws.Range("P12").Value = 1
$ws.Range("R12").Value = 'Knärot'
$ws.Range("S12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 24618-2022 artfynd.xlsx", "A 24618-2022")'
$ws.Range("T12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 24618-2022 karta.png", "A 24618-2022")'
$ws.Range("U12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/knärot/A 24618-2022 karta knärot.png", "A 24618-2022")'
$ws.Range("V12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 24618-2022 FSC-klagomål.docx", "A 24618-2022")'
$ws.Range("W12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 24618-2022 FSC-klagomål mail.docx", "A 24618-2022")'
$ws.Range("X12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 24618-2022 tillsynsbegäran.docx", "A 24618-2022")'
$ws.Range("Y12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 24618-2022 tillsynsbegäran mail.docx", "A 24618-2022")'
# Row 13
$ws.Range("C13").Value = 46077
# Row 14
$ws.Range("C14").Value = 46077
# Row 15
$ws.Range("A15").Value = 'A 56799-2022'
$ws.Range("B15").Value = 44894.425625
$ws.Range("C15").Value = 46077
$ws.Range("F15").ClearContents()
$ws.Range("G15").Value = 0.4
# Row 16
$ws.Range("A16").Value = 'A 51008-2023'
$ws.Range("B16").Value = 45218
$ws.Range("C16").Value = 46077
$ws.Range("F16").Value = 'Kommuner'
$ws.Range("G16").Value = 0.5
# Row 17
$ws.Range("A17").Value = 'A 57893-2023'
$ws.Range("B17").Value = 45247
$ws.Range("C17").Value = 46077
$ws.Range("F17").Value = 'Kommuner'
$ws.Range("G17").Value = 0.7
# Row 18
$ws.Range("A18").Value = 'A 37934-2024'
$ws.Range("B18").Value = 45544
$ws.Range("C18").Value = 46077
$ws.Range("G18").Value = 1.3
# Row 19
$ws.Range("A19").Value = 'A 55069-2023'
$ws.Range("C19").Value = 46077
$ws.Range("G19").Value = 3.2
# Row 20
$ws.Range("A20").Value = 'A 41895-2023'
$ws.Range("B20").Value = 45176
$ws.Range("C20").Value = 46077
$ws.Range("F20").Value = 'Övriga Aktiebolag'
$ws.Range("G20").Value = 0.6
# Row 21
$ws.Range("A21").Value = 'A 41899-2023'
$ws.Range("B21").Value = 45176
$ws.Range("C21").Value = 46077
$ws.Range("F21").Value = 'Övriga Aktiebolag'
$ws.Range("G21").Value = 1.8
# Row 22
$ws.Range("A22").Value = 'A 58383-2024'
$ws.Range("B22").Value = 45632
$ws.Range("C22").Value = 46077
$ws.Range("G22").Value = 1.6
# Row 23
$ws.Range("A23").Value = 'A 53369-2024'
$ws.Range("B23").Value = 45614
$ws.Range("C23").Value = 46077
$ws.Range("F23").ClearContents()
$ws.Range("G23").Value = 3.4
# Row 24
$ws.Range("A24").Value = 'A 13356-2022'
$ws.Range("B24").Value = 44645
$ws.Range("C24").Value = 46077
$ws.Range("F24").ClearContents()
$ws.Range("G24").Value = 1.6
# Row 25
$ws.Range("A25").Value = 'A 42951-2023'
$ws.Range("B25").Value = 45182
$ws.Range("C25").Value = 46077
$ws.Range("F25").Value = 'Övriga Aktiebolag'
$ws.Range("G25").Value = 0.8
# Row 26
$ws.Range("A26").Value = 'A 42955-2023'
$ws.Range("B26").Value = 45182
$ws.Range("C26").Value = 46077
$ws.Range("F26").Value = 'Övriga Aktiebolag'
$ws.Range("G26").Value = 0.3
# Row 27
$ws.Range("A27").Value = 'A 42960-2023'
$ws.Range("B27").Value = 45182
$ws.Range("C27").Value = 46077
$ws.Range("G27").Value = 0.4
# Row 28
$ws.Range("C28").Value = 46077
# Row 29
$ws.Range("A29").Value = 'A 55066-2023'
$ws.Range("B29").Value = 45237
$ws.Range("C29").Value = 46077
$ws.Range("F29").Value = 'Kommuner'
$ws.Range("G29").Value = 0.9
# Row 30
$ws.Range("A30").Value = 'A 45423-2025'
$ws.Range("B30").Value = 45922.45137731481
$ws.Range("C30").Value = 46077
$ws.Range("G30").Value = 1.1
# Row 31
$ws.Range("A31").Value = 'A 2769-2023'
$ws.Range("B31").Value = 44944
$ws.Range("C31").Value = 46077
$ws.Range("G31").Value = 8.699999999999999
# Row 32
$ws.Range("A32").Value = 'A 41661-2025'
$ws.Range("B32").Value = 45902
$ws.Range("C32").Value = 46077
$ws.Range("F32").ClearContents()
$ws.Range("G32").Value = 5.8
# Row 33
$ws.Range("A33").Value = 'A 42957-2023'
$ws.Range("B33").Value = 45182
$ws.Range("C33").Value = 46077
$ws.Range("F33").Value = 'Övriga Aktiebolag'
$ws.Range("G33").Value = 0.5
